$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 105
$ws.Range("G3").Value = 114
$ws.Range("I3").Value = 70
$ws.Range("G4").Value = 13
$ws.Range("G5").Value = 13
$ws.Range("G6").Value = 14
$ws.Range("I6").Value = 78
$ws.Range("G7").Value = 73
$ws.Range("I7").Value = 19
$ws.Range("G8").Value = 79
$ws.Range("I8").Value = 13
$ws.Range("G9").Value = 2
$ws.Range("I9").Value = 90
$ws.Range("G10").Value = 173
$ws.Range("I10").Value = 11
$ws.Range("I11").Value = 51
$ws.Range("G12").Value = 40
$ws.Range("I12").Value = 142
$ws.Range("G13").Value = 16
$ws.Range("I13").Value = 167
$ws.Range("G14").Value = 16
$ws.Range("I14").Value = 167
$ws.Range("G15").Value = 129
$ws.Range("I15").Value = 54
$ws.Range("G16").Value = 49
$ws.Range("I16").Value = 43
$ws.Range("G17").Value = 49
$ws.Range("I17").Value = 42
$ws.Range("G18").Value = 11
$ws.Range("I18").Value = 81
$ws.Range("G19").Value = 6
$ws.Range("I19").Value = 86
$ws.Range("G20").Value = 52
$ws.Range("I20").Value = 39
$ws.Range("G21").Value = 48
$ws.Range("I21").Value = 43
$ws.Range("G22").Value = 48
$ws.Range("I22").Value = 43
$ws.Range("I23").Value = 10
$ws.Range("I24").Value = 10
$ws.Range("G25").Value = 31
$ws.Range("G26").Value = 75
$ws.Range("I26").Value = 109
$ws.Range("G27").Value = 117
$ws.Range("I27").Value = 67
$ws.Range("G28").Value = 150
$ws.Range("I28").Value = 33
$ws.Range("G29").Value = 52
$ws.Range("I29").Value = 39
$ws.Range("G30").Value = 11
$ws.Range("I30").Value = 81
$ws.Range("G31").Value = 82
$ws.Range("I31").Value = 10
$ws.Range("G32").Value = 80
$ws.Range("I32").Value = 12
$ws.Range("I33").Value = 20
$ws.Range("G34").Value = 59
$ws.Range("I34").Value = 32
$ws.Range("G35").Value = 13
$ws.Range("I35").Value = 79
$ws.Range("G36").Value = 28
$ws.Range("I36").Value = 64
$ws.Range("G37").Value = 46
$ws.Range("I37").Value = 45
$ws.Range("G38").Value = 128
$ws.Range("I38").Value = 55
$ws.Range("G39").Value = 172
$ws.Range("I39").Value = 12
$ws.Range("G40").Value = 58
$ws.Range("I40").Value = 124
$ws.Range("G41").Value = 128
$ws.Range("I41").Value = 55
$ws.Range("I42").Value = 30
$ws.Range("G43").Value = 137
$ws.Range("I43").Value = 46
$ws.Range("G44").Value = 136
$ws.Range("I44").Value = 47
$ws.Range("G45").Value = 129
$ws.Range("I45").Value = 54
$ws.Range("G46").Value = 129
$ws.Range("I46").Value = 54
$ws.Range("G47").Value = 115
$ws.Range("I47").Value = 69
$ws.Range("G48").Value = 117
$ws.Range("I48").Value = 67
$ws.Range("G49").Value = 83
$ws.Range("I49").Value = 101
$ws.Range("G50").Value = 10
$ws.Range("I50").Value = 173
$ws.Range("G51").Value = 113
$ws.Range("I51").Value = 71
$ws.Range("G52").Value = 96
$ws.Range("I52").Value = 88
$ws.Range("G53").Value = 136
$ws.Range("I53").Value = 47
$ws.Range("G54").Value = 39
$ws.Range("I54").Value = 52
$ws.Range("G55").Value = 8
$ws.Range("G56").Value = 84
$ws.Range("I56").Value = 8
$ws.Range("G57").Value = 39
$ws.Range("I57").Value = 52
$ws.Range("G58").Value = 166
$ws.Range("I58").Value = 18
$ws.Range("G59").Value = 44
$ws.Range("I59").Value = 47
$ws.Range("G60").Value = 30
$ws.Range("G61").Value = 16
$ws.Range("I61").Value = 78
$ws.Range("G62").Value = 177
$ws.Range("I62").Value = 7
$ws.Range("G63").Value = 41
$ws.Range("I63").Value = 141
$ws.Range("G64").Value = 27
$ws.Range("G65").Value = 27
$ws.Range("G66").Value = 146
$ws.Range("I66").Value = 37
$ws.Range("G67").Value = 57
$ws.Range("I67").Value = 34
$ws.Range("G68").Value = 39
$ws.Range("I68").Value = 52
$ws.Range("G69").Value = 86
$ws.Range("I69").Value = 6
$ws.Range("G70").Value = 30
$ws.Range("I70").Value = 61
$ws.Range("G71").Value = 30
$ws.Range("G72").Value = 28
$ws.Range("I72").Value = 64
$ws.Range("G73").Value = 46
$ws.Range("I73").Value = 45
$ws.Range("G74").Value = 14
$ws.Range("I74").Value = 169
$ws.Range("G75").Value = 118
$ws.Range("I75").Value = 66
$ws.Range("G76").Value = 75
$ws.Range("I76").Value = 109
$ws.Range("G77").Value = 49
$ws.Range("I77").Value = 133
$ws.Range("G78").Value = 9
$ws.Range("I78").Value = 174
$ws.Range("G79").Value = 156
$ws.Range("I79").Value = 27
$ws.Range("I80").Value = 19
$ws.Range("G81").Value = 131
$ws.Range("I81").Value = 52
$ws.Range("I82").Value = 53
$ws.Range("G83").Value = 60
$ws.Range("I83").Value = 122
$ws.Range("G84").Value = 181
$ws.Range("I84").Value = 3
$ws.Range("G85").Value = 179
$ws.Range("I85").Value = 5
$ws.Range("G86").Value = 150
$ws.Range("I86").Value = 33
$ws.Range("G87").Value = 41
$ws.Range("I87").Value = 325
$ws.Range("G88").Value = 132
$ws.Range("I88").Value = 51
$ws.Range("G89").Value = 132
$ws.Range("I89").Value = 51
$ws.Range("G90").Value = 43
$ws.Range("I90").Value = 139
$ws.Range("G91").Value = 135
$ws.Range("I91").Value = 52
$ws.Range("G92").Value = 48
$ws.Range("I92").Value = 134
$ws.Range("G93").Value = 131
$ws.Range("I93").Value = 52
$ws.Range("G94").Value = 131
$ws.Range("I94").Value = 52
$ws.Range("G95").Value = 2
$ws.Range("I95").Value = 180
$ws.Range("G96").Value = 113
$ws.Range("I96").Value = 71
$ws.Range("G97").Value = 94
$ws.Range("I97").Value = 90
$ws.Range("G98").Value = 41
$ws.Range("I98").Value = 141
$ws.Range("G99").Value = 142
$ws.Range("I99").Value = 41
$ws.Range("G100").Value = 112
$ws.Range("I100").Value = 72
$ws.Range("G101").Value = 37
$ws.Range("I101").Value = 145
$ws.Range("G102").Value = 142
$ws.Range("I102").Value = 41
$ws.Range("G103").Value = 113
$ws.Range("I103").Value = 71
$ws.Range("G104").Value = 113
$ws.Range("I104").Value = 71
$ws.Range("F105").Value = 45234
$ws.Range("G105").Value = 1
$ws.Range("G106").Value = 150
$ws.Range("I106").Value = 33
$ws.Range("G107").Value = 50
$ws.Range("I107").Value = 41
$ws.Range("G108").Value = 50
$ws.Range("I108").Value = 41
$ws.Range("G109").Value = 50
$ws.Range("I109").Value = 41
$ws.Range("G110").Value = 87
$ws.Range("I110").Value = 5
$ws.Range("G111").Value = 24
$ws.Range("I111").Value = 68
$ws.Range("G112").Value = 63
$ws.Range("I112").Value = 28
$ws.Range("G113").Value = 29
$ws.Range("I113").Value = 154
$ws.Range("G114").Value = 32
$ws.Range("I114").Value = 60
$ws.Range("G115").Value = 17
$ws.Range("I115").Value = 75
$ws.Range("G116").Value = 74
$ws.Range("I116").Value = 18
$ws.Range("G117").Value = 120
$ws.Range("I117").Value = 64
$ws.Range("G118").Value = 141
$ws.Range("I118").Value = 42
$ws.Range("G119").Value = 131
$ws.Range("I119").Value = 52
$ws.Range("G120").Value = 131
$ws.Range("I120").Value = 52
$ws.Range("G121").Value = 131
$ws.Range("I121").Value = 52
$ws.Range("G122").Value = 131
$ws.Range("I122").Value = 52
$ws.Range("G123").Value = 36
$ws.Range("I123").Value = 56
$ws.Range("G124").Value = 38
$ws.Range("I124").Value = 144
$ws.Range("G125").Value = 22
$ws.Range("I125").Value = 70
$ws.Range("G126").Value = 9
$ws.Range("I126").Value = 83
$ws.Range("G127").Value = 23
$ws.Range("I127").Value = 69
$ws.Range("G128").Value = 17
$ws.Range("I128").Value = 75
$ws.Range("G129").Value = 77
$ws.Range("I129").Value = 15
$ws.Range("G130").Value = 64
$ws.Range("I130").Value = 27
$ws.Range("G131").Value = 32
$ws.Range("I131").Value = 60
$ws.Range("G132").Value = 91
$ws.Range("I132").Value = 1
$ws.Range("G133").Value = 79
$ws.Range("I133").Value = 13
$ws.Range("G134").Value = 89
$ws.Range("I134").Value = 3
$ws.Range("G135").Value = 3
$ws.Range("I135").Value = 89
$ws.Range("G136").Value = 39
$ws.Range("I136").Value = 52
$ws.Range("G137").Value = 38
$ws.Range("I137").Value = 53
$ws.Range("G138").Value = 56
$ws.Range("I138").Value = 35
$ws.Range("G139").Value = 38
$ws.Range("I139").Value = 53
$ws.Range("G140").Value = 53
$ws.Range("I140").Value = 38
$ws.Range("G141").Value = 27
$ws.Range("I141").Value = 65
$ws.Range("G142").Value = 85
$ws.Range("I142").Value = 7
$ws.Range("I143").Value = 3
$ws.Range("G144").Value = 64
$ws.Range("I144").Value = 27
$ws.Range("G145").Value = 41
$ws.Range("I145").Value = 50
$ws.Range("G146").Value = 24
$ws.Range("G147").Value = 70
$ws.Range("I147").Value = 22
$ws.Range("G148").Value = 69
$ws.Range("I148").Value = 23
$ws.Range("G149").Value = 50
$ws.Range("I149").Value = 41
$ws.Range("G150").Value = 17
$ws.Range("I150").Value = 166
$ws.Range("G151").Value = 61
$ws.Range("I151").Value = 121
$ws.Range("G152").Value = 146
$ws.Range("I152").Value = 37
$ws.Range("G153").Value = 111
$ws.Range("I153").Value = 73
$ws.Range("G154").Value = 111
$ws.Range("I154").Value = 73
$ws.Range("G155").Value = 26
$ws.Range("I155").Value = 157
$ws.Range("G156").Value = 26
$ws.Range("I156").Value = 157
$ws.Range("G157").Value = 26
$ws.Range("I157").Value = 157
$ws.Range("G158").Value = 26
$ws.Range("I158").Value = 157
$ws.Range("G159").Value = 79
$ws.Range("I159").Value = 12
$ws.Range("G160").Value = 82
$ws.Range("I160").Value = 10
$ws.Range("G161").Value = 100
$ws.Range("I161").Value = 84
$ws.Range("G162").Value = 144
$ws.Range("I162").Value = 39
$ws.Range("G163").Value = 144
$ws.Range("I163").Value = 39
$ws.Range("G164").Value = 97
$ws.Range("I164").Value = 87
$ws.Range("G165").Value = 97
$ws.Range("I165").Value = 87
$ws.Range("G166").Value = 150
$ws.Range("I166").Value = 33
$ws.Range("G167").Value = 150
$ws.Range("I167").Value = 33
$ws.Range("I168").Value = 91
$ws.Range("G169").Value = 62
$ws.Range("I169").Value = 120
$ws.Range("G170").Value = 179
$ws.Range("I170").Value = 5
$ws.Range("I171").Value = 3
$ws.Range("G172").Value = 94
$ws.Range("I172").Value = 90
$ws.Range("G173").Value = 15
$ws.Range("I173").Value = 168
$ws.Range("G174").Value = 64
$ws.Range("I174").Value = 27
$ws.Range("G175").Value = 38
$ws.Range("I175").Value = 53
$ws.Range("G176").Value = 73
$ws.Range("I176").Value = 19
$ws.Range("G177").Value = 174
$ws.Range("I177").Value = 10
$ws.Range("G178").Value = 23
$ws.Range("I178").Value = 69
$ws.Range("G179").Value = 44
$ws.Range("I179").Value = 47
$ws.Range("G180").Value = 147
$ws.Range("I180").Value = 36
$ws.Range("G181").Value = 43
$ws.Range("I181").Value = 48
$ws.Range("G182").Value = 12
$ws.Range("I182").Value = 80
$ws.Range("G183").Value = 7
$ws.Range("I183").Value = 85
$ws.Range("G184").Value = 85
$ws.Range("I184").Value = 7
$ws.Range("G185").Value = 25
$ws.Range("G186").Value = 14
$ws.Range("I186").Value = 78
$ws.Range("G187").Value = 46
$ws.Range("I187").Value = 45
$ws.Range("G188").Value = 42
$ws.Range("I188").Value = 49
$ws.Range("G189").Value = 70
$ws.Range("I189").Value = 22
$ws.Range("G190").Value = 33
$ws.Range("I190").Value = 150
$ws.Range("G191").Value = 154
$ws.Range("I191").Value = 29
$ws.Range("G192").Value = 154
$ws.Range("I192").Value = 29
$ws.Range("I193").Value = 84
$ws.Range("I194").Value = 84
$ws.Range("G195").Value = 5
$ws.Range("I195").Value = 177
$ws.Range("G196").Value = 5
$ws.Range("I196").Value = 177
$ws.Range("G197").Value = 21
$ws.Range("I197").Value = 162
$ws.Range("G198").Value = 21
$ws.Range("I198").Value = 162
$ws.Range("G199").Value = 21
$ws.Range("I199").Value = 162
$ws.Range("G200").Value = 21
$ws.Range("I200").Value = 162
$ws.Range("G201").Value = 168
$ws.Range("I201").Value = 16
$ws.Range("G202").Value = 168
$ws.Range("I202").Value = 16
$ws.Range("G203").Value = 168
$ws.Range("I203").Value = 16
$ws.Range("G204").Value = 168
$ws.Range("I204").Value = 16
$ws.Range("G205").Value = 57
$ws.Range("I205").Value = 34
$ws.Range("I206").Value = 17
$ws.Range("G207").Value = 148
$ws.Range("I207").Value = 35
$ws.Range("G208").Value = 136
$ws.Range("I208").Value = 47
$ws.Range("G209").Value = 60
$ws.Range("I209").Value = 121
$ws.Range("G210").Value = 60
$ws.Range("I210").Value = 121
$ws.Range("G211").Value = 69
$ws.Range("I211").Value = 115
$ws.Range("G212").Value = 69
$ws.Range("I212").Value = 115
$ws.Range("G213").Value = 145
$ws.Range("I213").Value = 38
$ws.Range("G214").Value = 129
$ws.Range("I214").Value = 54
$ws.Range("G215").Value = 9
$ws.Range("I215").Value = 83
$ws.Range("G216").Value = 11
$ws.Range("I216").Value = 81
$ws.Range("G217").Value = 85
$ws.Range("I217").Value = 7
$ws.Range("G218").Value = 48
$ws.Range("I218").Value = 43
$ws.Range("G219").Value = 53
$ws.Range("I219").Value = 38
$ws.Range("G220").Value = 39
$ws.Range("I220").Value = 52
$ws.Range("G221").Value = 35
$ws.Range("I221").Value = 57
$ws.Range("G222").Value = 8
$ws.Range("I222").Value = 84
$ws.Range("G223").Value = 87
$ws.Range("I223").Value = 5
$ws.Range("G224").Value = 43
$ws.Range("I224").Value = 48
$ws.Range("G225").Value = 99
$ws.Range("I225").Value = 85
$ws.Range("G226").Value = 99
$ws.Range("I226").Value = 85
$ws.Range("G227").Value = 141
$ws.Range("I227").Value = 42
$ws.Range("G228").Value = 154
$ws.Range("I228").Value = 29
$ws.Range("G229").Value = 136
$ws.Range("I229").Value = 47
$ws.Range("G230").Value = 34
$ws.Range("I230").Value = 149
$ws.Range("G231").Value = 21
$ws.Range("I231").Value = 162
$ws.Range("I232").Value = 59
$ws.Range("G233").Value = 156
$ws.Range("I233").Value = 27
$ws.Range("G234").Value = 36
$ws.Range("I234").Value = 147
$ws.Range("G235").Value = 31
$ws.Range("I235").Value = 152
$ws.Range("G236").Value = 31
$ws.Range("I236").Value = 152
$ws.Range("G237").Value = 155
$ws.Range("I237").Value = 28
$ws.Range("G238").Value = 6
$ws.Range("I238").Value = 177
$ws.Range("G239").Value = 31
$ws.Range("I239").Value = 61
$ws.Range("G240").Value = 31
$ws.Range("I240").Value = 61
$ws.Range("G241").Value = 36
$ws.Range("I241").Value = 147
$ws.Range("G242").Value = 165
$ws.Range("I242").Value = 19
$ws.Range("G243").Value = 67
$ws.Range("I243").Value = 25
$ws.Range("G244").Value = 25
$ws.Range("I244").Value = 158
$ws.Range("G245").Value = 158
$ws.Range("I245").Value = 25
$ws.Range("G246").Value = 135
$ws.Range("I246").Value = 48
$ws.Range("I247").Value = 46
$ws.Range("G248").Value = 8
$ws.Range("G249").Value = 134
$ws.Range("I249").Value = 48
$ws.Range("G250").Value = 9
$ws.Range("G251").Value = 77
$ws.Range("I251").Value = 15
$ws.Range("G252").Value = 35
$ws.Range("I252").Value = 57
$ws.Range("I253").Value = 82
$ws.Range("G254").Value = 13
$ws.Range("I254").Value = 170
$ws.Range("G255").Value = 14
$ws.Range("I255").Value = 78
$ws.Range("G256").Value = 78
$ws.Range("I256").Value = 106
$ws.Range("G257").Value = 172
$ws.Range("I257").Value = 12
$ws.Range("I258").Value = 19
$ws.Range("I259").Value = 251
$ws.Range("G260").Value = 162
$ws.Range("I260").Value = 204
$ws.Range("I261").Value = 88
$ws.Range("I262").Value = 68

Write-Output "Applied 489 cell updates"